$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the task description for row 5 (shorten it)
$ws.Range("D5").Value = "Implement backspace, more functionality"

# Mark row 5 as Done (Yes) - this also triggers the "Good" (Gut) cell style
$ws.Range("G5").Value = "Yes"
$ws.Range("G5").Style = "Gut"

# Move the active selection from D17 to D16
$ws.Range("D16").Select()
